# Update: po 17. 01. 2022
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Corrections to existing rows (AgTests / AgPosit columns) ---
$ws.Range("F558").Value = 25169

$ws.Range("F656").Value = 52373
$ws.Range("F657").Value = 34054

$ws.Range("F664").Value = 26408

$ws.Range("F665").Value = 28212
$ws.Range("G665").Value = 648

$ws.Range("F667").Value = 16803
$ws.Range("G667").Value = 598

$ws.Range("F669").Value = 23388

$ws.Range("F670").Value = 52382
$ws.Range("G670").Value = 904

$ws.Range("F671").Value = 32590
$ws.Range("G671").Value = 614

$ws.Range("F672").Value = 29725
$ws.Range("G672").Value = 579

$ws.Range("F673").Value = 10078
$ws.Range("G673").Value = 308

$ws.Range("F674").Value = 28608
$ws.Range("G674").Value = 680

$ws.Range("F675").Value = 13434
$ws.Range("G675").Value = 360

$ws.Range("F676").Value = 28016
$ws.Range("G676").Value = 445

$ws.Range("F677").Value = 55875
$ws.Range("G677").Value = 796

$ws.Range("F678").Value = 33647
$ws.Range("G678").Value = 522

$ws.Range("F679").Value = 29216
$ws.Range("G679").Value = 516

$ws.Range("F680").Value = 27791
$ws.Range("G680").Value = 529

# --- New rows 681-683 ---
$ws.Range("A681").Value = 44575
$ws.Range("A681").NumberFormat = "yyyy-mm-dd"
$ws.Range("B681").Value = 876448
$ws.Range("C681").Value = 9013
$ws.Range("D681").Value = 2049
$ws.Range("E681").Value = 17252
$ws.Range("F681").Value = 22432
$ws.Range("G681").Value = 490

$ws.Range("A682").Value = 44576
$ws.Range("A682").NumberFormat = "yyyy-mm-dd"
$ws.Range("B682").Value = 879656
$ws.Range("C682").Value = 12981
$ws.Range("D682").Value = 3208
$ws.Range("E682").Value = 17300
$ws.Range("F682").Value = 10193
$ws.Range("G682").Value = 357

$ws.Range("A683").Value = 44577
$ws.Range("A683").NumberFormat = "yyyy-mm-dd"
$ws.Range("B683").Value = 880671
$ws.Range("C683").Value = 4940
$ws.Range("D683").Value = 1015
$ws.Range("E683").Value = 17352
$ws.Range("F683").Value = 14017
$ws.Range("G683").Value = 431
